$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column J with header info and one data value for row 10 (CAOCAO description column)
$ws.Range("J1").Value = "desc"
$ws.Range("J2").Value = "string"
$ws.Range("J3").Value = "描述"
$ws.Range("J10").Value = "CAOCAO_DESC"

# Update active selection as seen in the edited workbook
$ws.Range("E10").Select()
